$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: H2 ("productId") no longer holds a value - became a blank (foreign-key)
# cell, keeping its existing style.
$ws.Range("H2").Value = ""

# Row 12: A12 ("cart") and B12 ("cartId") are cleared out (the "cart" table /
# cartId column were dropped from the model), keeping their existing styles.
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""

# Row 12: H12 used to reference "cartId" - now references "userId" instead,
# so copy the visual style used for the other "userId" cells (e.g. H3) over
# before writing the new text.
$ws.Range("H3").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = "userId"

# Row 13: B13 ("uid") column is removed entirely from the model.
$ws.Range("B13").ClearContents()

# Update the active selection to match the new cursor position.
$ws.Range("G10").Select()
